# Add a "filtercars" worksheet after "carbrandtest" with a
# carbrand / carTitle lookup table used for budget-based search tests.

$wb = $excel.ActiveWorkbook

$wsBrand = $wb.Worksheets.Item("carbrandtest")

# carbrandtest stops being the active/selected tab; its selection becomes
# a full first-row selection (no specific active cell focus).
$wsBrand.Activate()
$wsBrand.Rows.Item(1).Select()

# Insert the new sheet right after carbrandtest so it lands as the 3rd tab.
$wsFilter = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $wsBrand)
$wsFilter.Name = "filtercars"

# Header row
$wsFilter.Range("A1").Value = "carbrand"
$wsFilter.Range("B1").Value = "carTitle"

# Data rows
$wsFilter.Range("A2").Value = "Honda"
$wsFilter.Range("B2").Value = "Honda Cars"

$wsFilter.Range("A3").Value = "Toyota"
$wsFilter.Range("B3").Value = "Toyota Cars"

$wsFilter.Range("A4").Value = "BMW"
$wsFilter.Range("B4").Value = "BMW Cars"

$wsFilter.Range("A5").Value = "Hyundai"
$wsFilter.Range("B5").Value = "Hyundai Cars"

# Column widths (23.08984375 / 15 characters)
$wsFilter.Columns.Item(1).ColumnWidth = 22.333333333333336
$wsFilter.Columns.Item(2).ColumnWidth = 14.166666666666666

# Final selection/active cell on the new sheet
$wsFilter.Range("D5").Select()
